# Add scripts for validation data feature
# Replace a few identifier values on the data row (row 2) with new
# validation-data style identifiers, and populate the previously blank
# "Quote Start Date" / "Quote End Date" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2 ("Exp Doc Nbr") and AM2 ("Func Loc") already carry a text number
# format, so plain assignment keeps them stored as text.
$ws.Range("F2").Value = "215243005501"
$ws.Range("AM2").Value = "ABCD228926"

# AG2 ("SAID") is formatted as a plain integer, so round-trip through a
# text format to force the new value to be stored as text (it looks
# numeric but must stay a string), then restore the original numeric
# display format.
$ws.Range("AG2").NumberFormat = "@"
$ws.Range("AG2").Value = "103373890322"
$ws.Range("AG2").NumberFormat = "0"

# DD2 ("Quote Start Date") / DE2 ("Quote End Date") were empty; fill them
# in using the same date format as the neighboring "Close Date" column.
$ws.Range("DD2").NumberFormat = "@"
$ws.Range("DD2").Value = "04/13/2016"
$ws.Range("DD2").NumberFormat = "m/d/yy"
$ws.Range("DE2").NumberFormat = "m/d/yy"
$ws.Range("DE2").Value = 42595

# Excel resets the reference style to A1 (clearing the R1C1 calc hint)
# and the view scrolls back to the top with a fresh selection.
$excel.ReferenceStyle = -4150
$ws.Range("F7").Select()

$wb.Save()
